$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "請求書"

$ws.Range("B4").Value = "株式会社ABC"
$ws.Range("F4").Value = "No."
$ws.Range("G4").Value = "'0001"

$ws.Range("B5").Value = "〒101-0022 東京都千代田区神田練塀町300"
$ws.Range("F5").Value = "日付"
$ws.Range("G5").Value = "'2024/12/15"

$ws.Range("B6").Value = "TEL:03-1234-5678 FAX:03-1234-5678"

$ws.Range("B7").Value = "担当者名:鈴木一郎 様"

$ws.Range("B10").Value = "商品名"
$ws.Range("C10").Value = "数量"
$ws.Range("D10").Value = "単価"
$ws.Range("E10").Value = "金額"

$ws.Range("B11").Value = "商品A"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 10000
$ws.Range("E11").Value = 20000

$ws.Range("B12").Value = "商品B"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 15000
$ws.Range("E12").Value = 15000

$ws.Range("E13").Value = 35000

$ws.Range("B15").Value = "小計"
$ws.Range("E15").Value = 35000

$ws.Range("B16").Value = "消費税"
$ws.Range("E16").Value = 3500

$ws.Range("B17").Value = "合計"
$ws.Range("E17").Value = 38500
